# Apply scheduled-runner market/profit data refresh to each sheet.
# Mirrors the upstream OOXML diff: per-row cells for currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ), and LeveProfit(NQ/HQ) are updated to the latest cached values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 33278.53
$ws.Range("J17").Value = 33278.53
$ws.Range("L17").Value = 99835.59
$ws.Range("N17").Value = -100171.59
# Row 58
$ws.Range("H58").Value = 22673.877
$ws.Range("I58").Value = 350
$ws.Range("J58").Value = 24658.223
$ws.Range("K58").Value = 1050
$ws.Range("L58").Value = 73974.66900000001
$ws.Range("M58").Value = -900
$ws.Range("N58").Value = -74274.66900000001
# Row 74
$ws.Range("H74").Value = 3577.375
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3577.375
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3577.375
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -5449.375
# Row 76
$ws.Range("H76").Value = 2872.0625
$ws.Range("I76").Value = 2872.0625
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2872.0625
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2557.0625
$ws.Range("N76").ClearContents()
# Row 77
$ws.Range("H77").Value = 3577.375
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3577.375
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 17886.875
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -27246.875
# Row 79
$ws.Range("H79").Value = 2872.0625
$ws.Range("I79").Value = 2872.0625
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2872.0625
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1780.0625
$ws.Range("N79").ClearContents()
# Row 82
$ws.Range("H82").Value = 6188.4
# Row 85
$ws.Range("H85").Value = 6188.4
# Row 115
$ws.Range("H115").Value = 2444.4443
$ws.Range("J115").Value = 2444.4443
$ws.Range("L115").Value = 7333.3329
$ws.Range("N115").Value = -10467.3329
# Row 137
$ws.Range("H137").Value = 3500.8928
$ws.Range("I137").Value = 3521.875
$ws.Range("J137").Value = 3375
$ws.Range("K137").Value = 10565.625
$ws.Range("L137").Value = 10125
$ws.Range("M137").Value = -8015.625
$ws.Range("N137").Value = -15225
# Row 138
$ws.Range("H138").Value = 4214.8975
$ws.Range("I138").Value = 1316.6471
$ws.Range("J138").Value = 5022.6064
$ws.Range("K138").Value = 3949.9413
$ws.Range("L138").Value = 15067.8192
$ws.Range("M138").Value = 1190.0587
$ws.Range("N138").Value = -25347.8192

$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 2963.02
$ws.Range("I32").Value = 2303.5813
$ws.Range("J32").Value = 7013.857
$ws.Range("K32").Value = 2303.5813
$ws.Range("L32").Value = 7013.857
$ws.Range("M32").Value = -2016.5813
$ws.Range("N32").Value = -7587.857
# Row 61
$ws.Range("H61").Value = 3012.1292
$ws.Range("I61").Value = 1086
$ws.Range("J61").Value = 5066.6665
$ws.Range("K61").Value = 1086
$ws.Range("L61").Value = 5066.6665
$ws.Range("M61").Value = -874
$ws.Range("N61").Value = -5490.6665
# Row 110
$ws.Range("H110").Value = 1363.1052
$ws.Range("J110").Value = 2484
$ws.Range("L110").Value = 2484
$ws.Range("N110").Value = -6574
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 3012.1292
$ws.Range("I136").Value = 1086
$ws.Range("J136").Value = 5066.6665
$ws.Range("K136").Value = 3258
$ws.Range("L136").Value = 15199.9995
$ws.Range("M136").Value = -708
$ws.Range("N136").Value = -20299.9995

$ws = $wb.Worksheets.Item("BSM")

# Row 14
$ws.Range("H14").Value = 70009
$ws.Range("J14").Value = 70009
$ws.Range("L14").Value = 70009
$ws.Range("N14").Value = -70353
# Row 133
$ws.Range("H133").Value = 11582.223
$ws.Range("J133").Value = 11582.223
$ws.Range("L133").Value = 11582.223
$ws.Range("N133").Value = -21702.223
# Row 134
$ws.Range("H134").Value = 1972.8379
$ws.Range("I134").Value = 1441.0294
$ws.Range("K134").Value = 4323.0882
$ws.Range("M134").Value = -1788.0882

$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 2147.9
$ws.Range("I16").Value = 896.5
$ws.Range("J16").Value = 4025
$ws.Range("K16").Value = 896.5
$ws.Range("L16").Value = 4025
$ws.Range("M16").Value = -609.5
$ws.Range("N16").Value = -4599
# Row 31
$ws.Range("H31").Value = 4038.7693
$ws.Range("I31").Value = 2945.5652
$ws.Range("K31").Value = 2945.5652
$ws.Range("M31").Value = -2650.5652
# Row 34
$ws.Range("H34").Value = 4038.7693
$ws.Range("I34").Value = 2945.5652
$ws.Range("K34").Value = 2945.5652
$ws.Range("M34").Value = -2743.5652
# Row 58
$ws.Range("H58").Value = 9261636
$ws.Range("I58").Value = 1475.4706
$ws.Range("J58").Value = 25003910
$ws.Range("K58").Value = 1475.4706
$ws.Range("L58").Value = 25003910
$ws.Range("M58").Value = -1272.4706
$ws.Range("N58").Value = -25004316
# Row 113
$ws.Range("H113").Value = 2147.9
$ws.Range("I113").Value = 896.5
$ws.Range("J113").Value = 4025
$ws.Range("K113").Value = 896.5
$ws.Range("L113").Value = 4025
$ws.Range("M113").Value = 1273.5
$ws.Range("N113").Value = -8365
# Row 136
$ws.Range("H136").Value = 9261636
$ws.Range("I136").Value = 1475.4706
$ws.Range("J136").Value = 25003910
$ws.Range("K136").Value = 4426.4118
$ws.Range("L136").Value = 75011730
$ws.Range("M136").Value = -1876.4118
$ws.Range("N136").Value = -75016830

$ws = $wb.Worksheets.Item("CUL")

# Row 13
$ws.Range("H13").Value = 950
$ws.Range("I13").Value = 950
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 2850
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -2682
$ws.Range("N13").ClearContents()
# Row 47
$ws.Range("H47").Value = 2513
$ws.Range("J47").Value = 3284
$ws.Range("L47").Value = 9852
$ws.Range("N47").Value = -10714
# Row 120
$ws.Range("H120").Value = 19378.75
$ws.Range("I120").Value = 18757.5
$ws.Range("K120").Value = 56272.5
$ws.Range("M120").Value = -51434.5
# Row 131
$ws.Range("H131").Value = 1394.6818
$ws.Range("J131").Value = 1067.5333
$ws.Range("L131").Value = 3202.5999
$ws.Range("N131").Value = -13282.5999

$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 4377.7334
$ws.Range("I80").Value = 4432.5
$ws.Range("J80").Value = 4315.143
$ws.Range("K80").Value = 4432.5
$ws.Range("L80").Value = 4315.143
$ws.Range("M80").Value = -3434.5
$ws.Range("N80").Value = -6311.143
# Row 83
$ws.Range("H83").Value = 4377.7334
$ws.Range("I83").Value = 4432.5
$ws.Range("J83").Value = 4315.143
$ws.Range("K83").Value = 22162.5
$ws.Range("L83").Value = 21575.715
$ws.Range("M83").Value = -17170.5
$ws.Range("N83").Value = -31559.715
# Row 107
$ws.Range("H107").Value = 928.96875
$ws.Range("I107").Value = 446.35
$ws.Range("J107").Value = 1733.3334
$ws.Range("K107").Value = 446.35
$ws.Range("L107").Value = 1733.3334
$ws.Range("M107").Value = 1473.65
$ws.Range("N107").Value = -5573.3334
# Row 113
$ws.Range("H113").Value = 6005.5
$ws.Range("I113").Value = 2011
$ws.Range("K113").Value = 2011
$ws.Range("M113").Value = 159
# Row 132
$ws.Range("H132").Value = 3220.3171
$ws.Range("I132").Value = 2778.6667
$ws.Range("J132").Value = 4072.0715
$ws.Range("K132").Value = 8336.000100000001
$ws.Range("L132").Value = 12216.2145
$ws.Range("M132").Value = -5806.000100000001
$ws.Range("N132").Value = -17276.2145
# Row 137
$ws.Range("H137").Value = 29597
$ws.Range("J137").Value = 29597
$ws.Range("L137").Value = 29597
$ws.Range("N137").Value = -39797

$ws = $wb.Worksheets.Item("LTW")

# Row 82
$ws.Range("H82").Value = 2216.4285
$ws.Range("I82").Value = 1626
$ws.Range("J82").Value = 2544.4443
$ws.Range("K82").Value = 1626
$ws.Range("L82").Value = 2544.4443
$ws.Range("M82").Value = -1265
$ws.Range("N82").Value = -3266.4443
# Row 85
$ws.Range("H85").Value = 2216.4285
$ws.Range("I85").Value = 1626
$ws.Range("J85").Value = 2544.4443
$ws.Range("K85").Value = 1626
$ws.Range("L85").Value = 2544.4443
$ws.Range("M85").Value = -378
$ws.Range("N85").Value = -5040.4443
# Row 136
$ws.Range("H136").Value = 1450.2174
$ws.Range("I136").Value = 971.7143
$ws.Range("J136").Value = 2972.7273
$ws.Range("K136").Value = 2915.1429
$ws.Range("L136").Value = 8918.1819
$ws.Range("M136").Value = -365.1428999999998
$ws.Range("N136").Value = -14018.1819
# Row 138
$ws.Range("H138").Value = 29666.666
$ws.Range("J138").Value = 29666.666
$ws.Range("L138").Value = 29666.666
$ws.Range("N138").Value = -39946.666

$ws = $wb.Worksheets.Item("WVR")

# Row 46
$ws.Range("H46").Value = 43982.25
$ws.Range("J46").Value = 43982.25
$ws.Range("L46").Value = 43982.25
$ws.Range("N46").Value = -44444.25
# Row 62
$ws.Range("H62").Value = 5475.5
$ws.Range("I62").Value = 3951
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 3951
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -3327
$ws.Range("N62").Value = -8248
# Row 65
$ws.Range("H65").Value = 5475.5
$ws.Range("I65").Value = 3951
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 19755
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -16635
$ws.Range("N65").Value = -41240
# Row 118
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
# Row 120
$ws.Range("H120").Value = 22500
$ws.Range("J120").Value = 22500
$ws.Range("L120").Value = 22500
$ws.Range("N120").Value = -32176
# Row 132
$ws.Range("H132").Value = 46423.92
$ws.Range("I132").Value = 21102.666
$ws.Range("J132").Value = 54020.3
$ws.Range("K132").Value = 63307.99800000001
$ws.Range("L132").Value = 162060.9
$ws.Range("M132").Value = -60777.99800000001
$ws.Range("N132").Value = -167120.9
# Row 134
$ws.Range("H134").Value = 43982.25
$ws.Range("J134").Value = 43982.25
$ws.Range("L134").Value = 131946.75
$ws.Range("N134").Value = -137016.75
# Row 136
$ws.Range("H136").Value = 1134.8055
$ws.Range("I136").Value = 659.7308
$ws.Range("K136").Value = 1979.1924
$ws.Range("M136").Value = 570.8075999999999
